# Update cached market/profit figures (columns H-N) across the per-class
# leve sheets, as produced by the scheduled market-data refresh runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 156.36363
$ws.Range("I4").Value = 130.4762
$ws.Range("J4").Value = 700
$ws.Range("K4").Value = 130.4762
$ws.Range("L4").Value = 700
$ws.Range("M4").Value = -16.47620000000001
$ws.Range("N4").Value = -928

$ws.Range("H38").Value = 89.28570999999999
$ws.Range("I38").Value = 89.28570999999999
$ws.Range("K38").Value = 267.85713
$ws.Range("M38").Value = 104.14287

$ws.Range("H96").Value = 2828.2727
$ws.Range("I96").Value = 2467.25
$ws.Range("J96").Value = 3034.5715
$ws.Range("K96").Value = 7401.75
$ws.Range("L96").Value = 9103.7145
$ws.Range("M96").Value = -6028.75
$ws.Range("N96").Value = -11849.7145

$ws.Range("H99").Value = 3280.5715
$ws.Range("I99").Value = 988
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 2964
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = -1466
$ws.Range("N99").Value = -17996

$ws.Range("H132").Value = 25646396
$ws.Range("I132").Value = 40005984
$ws.Range("J132").Value = 4271.9287
$ws.Range("K132").Value = 120017952
$ws.Range("L132").Value = 12815.7861
$ws.Range("M132").Value = -120015422
$ws.Range("N132").Value = -17875.7861

$ws.Range("H137").Value = 3031.7778
$ws.Range("I137").Value = 2645.875
$ws.Range("J137").Value = 3340.5
$ws.Range("K137").Value = 7937.625
$ws.Range("L137").Value = 10021.5
$ws.Range("M137").Value = -5387.625
$ws.Range("N137").Value = -15121.5

$ws.Range("H138").Value = 5719.0625
$ws.Range("I138").Value = 2781.389
$ws.Range("J138").Value = 9496.071
$ws.Range("K138").Value = 8344.167000000001
$ws.Range("L138").Value = 28488.213
$ws.Range("M138").Value = -3204.167000000001
$ws.Range("N138").Value = -38768.213

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4191.943
$ws.Range("I61").Value = 3194
$ws.Range("J61").Value = 10179.6
$ws.Range("K61").Value = 3194
$ws.Range("L61").Value = 10179.6
$ws.Range("M61").Value = -2982
$ws.Range("N61").Value = -10603.6

$ws.Range("H74").Value = 1300.6316
$ws.Range("I74").Value = 1218.9412
$ws.Range("K74").Value = 1218.9412
$ws.Range("M74").Value = -344.9412

$ws.Range("H77").Value = 1300.6316
$ws.Range("I77").Value = 1218.9412
$ws.Range("K77").Value = 6094.706
$ws.Range("M77").Value = -1726.706

$ws.Range("H102").Value = 4656.7144
$ws.Range("I102").Value = 1994
$ws.Range("K102").Value = 1994
$ws.Range("M102").Value = -372

$ws.Range("H136").Value = 4191.943
$ws.Range("I136").Value = 3194
$ws.Range("J136").Value = 10179.6
$ws.Range("K136").Value = 9582
$ws.Range("L136").Value = 30538.8
$ws.Range("M136").Value = -7032
$ws.Range("N136").Value = -35638.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2081.348
$ws.Range("I20").Value = 1777.4286
$ws.Range("K20").Value = 1777.4286
$ws.Range("M20").Value = -1530.4286

$ws.Range("H82").Value = 31750.445
$ws.Range("J82").Value = 54191.4
$ws.Range("L82").Value = 54191.4
$ws.Range("N82").Value = -54957.4

$ws.Range("H85").Value = 31750.445
$ws.Range("J85").Value = 54191.4
$ws.Range("L85").Value = 54191.4
$ws.Range("N85").Value = -56843.4

$ws.Range("H122").Value = 78000
$ws.Range("J122").Value = 78000
$ws.Range("L122").Value = 78000
$ws.Range("N122").Value = -87800

$ws.Range("H134").Value = 2082.138
$ws.Range("I134").Value = 1817.64
$ws.Range("K134").Value = 5452.92
$ws.Range("M134").Value = -2917.92

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2092.8948
$ws.Range("I58").Value = 1886.7858
$ws.Range("J58").Value = 2670
$ws.Range("K58").Value = 1886.7858
$ws.Range("L58").Value = 2670
$ws.Range("M58").Value = -1683.7858
$ws.Range("N58").Value = -3076

$ws.Range("H132").Value = 1976.6428
$ws.Range("I132").Value = 1917.3
$ws.Range("K132").Value = 5751.9
$ws.Range("M132").Value = -3221.9

$ws.Range("H134").Value = 1160.6666
$ws.Range("I134").Value = 994.25
$ws.Range("J134").Value = 1826.3334
$ws.Range("K134").Value = 2982.75
$ws.Range("L134").Value = 5479.0002
$ws.Range("M134").Value = -447.75
$ws.Range("N134").Value = -10549.0002

$ws.Range("H136").Value = 2092.8948
$ws.Range("I136").Value = 1886.7858
$ws.Range("J136").Value = 2670
$ws.Range("K136").Value = 5660.357400000001
$ws.Range("L136").Value = 8010
$ws.Range("M136").Value = -3110.357400000001
$ws.Range("N136").Value = -13110

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 739.8946999999999
$ws.Range("J44").Value = 657.5454999999999
$ws.Range("L44").Value = 1972.6365
$ws.Range("N44").Value = -2768.6365

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9632.294
$ws.Range("I80").Value = 11206
$ws.Range("J80").Value = 6747.1665
$ws.Range("K80").Value = 11206
$ws.Range("L80").Value = 6747.1665
$ws.Range("M80").Value = -10208
$ws.Range("N80").Value = -8743.166499999999

$ws.Range("H83").Value = 9632.294
$ws.Range("I83").Value = 11206
$ws.Range("J83").Value = 6747.1665
$ws.Range("K83").Value = 56030
$ws.Range("L83").Value = 33735.8325
$ws.Range("M83").Value = -51038
$ws.Range("N83").Value = -43719.8325

$ws.Range("H92").Value = 16132.6
$ws.Range("J92").Value = 16132.6
$ws.Range("L92").Value = 16132.6
$ws.Range("N92").Value = -19876.6

$ws.Range("H97").Value = 2352.2083
$ws.Range("I97").Value = 1303.125
$ws.Range("J97").Value = 4450.375
$ws.Range("K97").Value = 1303.125
$ws.Range("L97").Value = 4450.375
$ws.Range("M97").Value = -807.125
$ws.Range("N97").Value = -5442.375

$ws.Range("H102").Value = 125796
$ws.Range("I102").Value = 1062
$ws.Range("J102").Value = 499998
$ws.Range("K102").Value = 1062
$ws.Range("L102").Value = 499998
$ws.Range("M102").Value = 560
$ws.Range("N102").Value = -503242

$ws.Range("H122").Value = 5282.846
$ws.Range("I122").Value = 3297.4
$ws.Range("J122").Value = 6523.75
$ws.Range("K122").Value = 9892.200000000001
$ws.Range("L122").Value = 19571.25
$ws.Range("M122").Value = -7442.200000000001
$ws.Range("N122").Value = -24471.25

$ws.Range("H132").Value = 3606.5974
$ws.Range("I132").Value = 3200.6167
$ws.Range("J132").Value = 5039.4707
$ws.Range("K132").Value = 9601.8501
$ws.Range("L132").Value = 15118.4121
$ws.Range("M132").Value = -7071.8501
$ws.Range("N132").Value = -20178.4121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8055.3335
$ws.Range("I7").Value = 4500
$ws.Range("K7").Value = 4500
$ws.Range("M7").Value = -4388

$ws.Range("H68").Value = 5122.222
$ws.Range("I68").Value = 5187.5
$ws.Range("K68").Value = 5187.5
$ws.Range("M68").Value = -4438.5

$ws.Range("H71").Value = 5122.222
$ws.Range("I71").Value = 5187.5
$ws.Range("K71").Value = 25937.5
$ws.Range("M71").Value = -22193.5

$ws.Range("H93").Value = 1658.28
$ws.Range("I93").Value = 1566.6818
$ws.Range("K93").Value = 1566.6818
$ws.Range("M93").Value = -318.6818000000001

$ws.Range("H122").Value = 15430.5
$ws.Range("I122").Value = 13541.842
$ws.Range("J122").Value = 20556.857
$ws.Range("K122").Value = 40625.526
$ws.Range("L122").Value = 61670.571
$ws.Range("M122").Value = -38175.526
$ws.Range("N122").Value = -66570.571

$ws.Range("H126").Value = 8055.3335
$ws.Range("I126").Value = 4500
$ws.Range("K126").Value = 13500
$ws.Range("M126").Value = -11030

$ws.Range("H136").Value = 4221.469
$ws.Range("I136").Value = 4261.9536
$ws.Range("K136").Value = 12785.8608
$ws.Range("M136").Value = -10235.8608

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 12157
$ws.Range("I76").Value = 12157
$ws.Range("K76").Value = 12157
$ws.Range("M76").Value = -11842

$ws.Range("H79").Value = 12157
$ws.Range("I79").Value = 12157
$ws.Range("K79").Value = 12157
$ws.Range("M79").Value = -11065

$ws.Range("H96").Value = 3004.1538
$ws.Range("I96").Value = 3058.125
$ws.Range("J96").Value = 2917.8
$ws.Range("K96").Value = 3058.125
$ws.Range("L96").Value = 2917.8
$ws.Range("M96").Value = -1685.125
$ws.Range("N96").Value = -5663.8

$ws.Range("H122").Value = 3747.5334
$ws.Range("J122").Value = 4934.7144
$ws.Range("L122").Value = 14804.1432
$ws.Range("N122").Value = -19704.1432

$ws.Range("H132").Value = 1000.7143
$ws.Range("I132").Value = 501
$ws.Range("K132").Value = 1503
$ws.Range("M132").Value = 1027
